# Applies the "prawie skonczona rejestracja" edit:
#  - "Menu i zakladki" sheet: split the single "check all controls" test
#    row into three rows (before login / as User / as Admin), expanding
#    the sheet's table accordingly.
#  - "Rejestracja" sheet: fix the duplicated step numbering in column A
#    (1,2,3,3,4,4,5,6 -> 1,2,3,4,5,6,7,8).
#  - Active tab moves from "Rejestracja" back to "Menu i zakladki", with
#    updated cell selections on both sheets.

$wb = $excel.ActiveWorkbook

$wsMenu = $wb.Worksheets.Item("Menu i zakładki")
$wsRejestracja = $wb.Worksheets.Item("Rejestracja")

# ---------------------------------------------------------------
# 1. "Menu i zakladki": reword row 2 and insert two new test rows.
# ---------------------------------------------------------------

$weryfikacja = $wsMenu.Range("C2").Value()
$odpowiedzialny = $wsMenu.Range("E2").Value()

$wsMenu.Range("B2").Value = "Sprawdzenie wszystkich kontolek (input, label, button) na wszystkich zakładkach przed zalogowaniem."
$wsMenu.Range("C2").Value = $weryfikacja
$wsMenu.Range("E2").Value = $odpowiedzialny

$wsMenu.Range("A3").Value = 2
$wsMenu.Range("B3").Value = "Sprawdzenie wszystkich kontolek (input, label, button) na wszystkich zakładkach jako User."
$wsMenu.Range("C3").Value = $weryfikacja
$wsMenu.Range("E3").Value = $odpowiedzialny

$wsMenu.Range("A4").Value = 3
$wsMenu.Range("B4").Value = "Sprawdzenie wszystkich kontolek (input, label, button) na wszystkich zakładkach jako Admin."
$wsMenu.Range("C4").Value = $weryfikacja
$wsMenu.Range("E4").Value = $odpowiedzialny

# Grow the worksheet table ("Tabela2") so it covers the two new rows.
$menuTable = $wsMenu.ListObjects.Item(1)
$null = $menuTable.Resize($wsMenu.Range("A1:E4"))

# Column B now holds longer text - widen it to fit.
$wsMenu.Columns.Item(2).ColumnWidth = 92.65

# ---------------------------------------------------------------
# 2. "Rejestracja": renumber column A sequentially (1..8).
# ---------------------------------------------------------------

$wsRejestracja.Range("A5").Value = 4
$wsRejestracja.Range("A6").Value = 5
$wsRejestracja.Range("A7").Value = 6
$wsRejestracja.Range("A8").Value = 7
$wsRejestracja.Range("A9").Value = 8

# ---------------------------------------------------------------
# 3. Move the active tab / selection back to "Menu i zakladki".
#    Do "Rejestracja" first so its selection updates without leaving
#    it as the tab that ends up active.
# ---------------------------------------------------------------

$null = $wsRejestracja.Range("B9").Select()
$null = $wsMenu.Range("B4").Select()

Write-Host "edit applied"
